# The workbook was re-uploaded with the single text cell edited and the
# selected cell moved. Reproduce both via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Este é um teste" -> "Este é o teste 2"
$ws.Range("A1").Value = "Este é o teste 2"

# Move/save the active selection from A2 to B3.
$ws.Range("B3").Select()
